$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new effort-row for Sara ("Fix requirements") right above the
# existing "Total effort" row (old row 25), pushing the total row down to 26.
$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = Get-Date -Year 2019 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Range("A25").NumberFormat = $ws.Range("A24").NumberFormat
$ws.Range("B25").Value = "Fix requirements"
$ws.Range("C25").Value = 2

# Copy the formatting of the row above so the new row matches the rest of
# the table (borders/fill/font).
$ws.Range("A24:C24").Copy()
$ws.Range("A25:C25").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Update the Total effort formula to include the newly inserted row.
$ws.Range("C26").Formula = "=SUM(C20:C25)"

# Update selection / view to match the target state.
$ws.Range("G22").Select()

# Match the target print setup (A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
